$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update EARNED (column C) values for core feature rows to 0
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("C9").Value = 0

# Update the view state of the sheet (scrolled position / selection)
$ws.Range("D13").Select()
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
